$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text content changes (case-insensitive delivery mechanism names) ---
$ws.Range("F2").Value = "deposit to card"
$ws.Range("F3").Value = "Referral"

# --- Extend sheet with a bordered block below the existing table (rows 4-10) ---
$boxRange = $ws.Range("A4:K10")
$boxRange.Interior.Color = 0xFFFFFF

$top = $ws.Range("A4:K4").Borders.Item(8)
$top.Color = 0xAAAAAA
$top.LineStyle = 1

$topAccentF = $ws.Range("F4").Borders.Item(8)
$topAccentF.Color = 0x999999
$topAccentF.LineStyle = 1

$topAccentI = $ws.Range("I4").Borders.Item(8)
$topAccentI.Color = 0x999999
$topAccentI.LineStyle = 1

$left = $ws.Range("A4:A10").Borders.Item(7)
$left.Color = 0xAAAAAA
$left.LineStyle = 1

$right = $ws.Range("K4:K10").Borders.Item(10)
$right.Color = 0xAAAAAA
$right.LineStyle = 1

$bottom = $ws.Range("A10:K10").Borders.Item(9)
$bottom.Color = 0xAAAAAA
$bottom.LineStyle = 1

for ($r = 4; $r -le 10; $r++) {
    $ws.Rows.Item($r).RowHeight = 13.55
}
